$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change swaps the report rows for the two handed-back files
# (b5a67361-...md and f84312db-...md) across all three sheets, and marks the
# b5a67361 file as freshly "Ready for handoff" with new handoff timestamps
# (matching a regenerated localization-status report).
# ---------------------------------------------------------------------------

# ---------------- Overview sheet ----------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 2 now describes f84312db (unchanged status/date), row 3 now describes
# b5a67361 with a fresh "Ready for handoff" status/date.
$ws1.Hyperlinks.Item(1).TextToDisplay = "f84312db-1505-4f53-a880-5d7ff6201955.md"
$ws1.Hyperlinks.Item(2).TextToDisplay = "b5a67361-cacf-47e3-aedd-3e6e3b8570f1.md"

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-23 12:52:49"

# ---------------- zh-cn sheet ----------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Hyperlinks 1-4 (A2,D2,F2,G2) display f84312db's links; 5-8 (A3,D3,F3,G3)
# display b5a67361's links (the underlying target URLs stay where they are).
$ws2.Hyperlinks.Item(1).TextToDisplay = "f84312db-1505-4f53-a880-5d7ff6201955.md"
$ws2.Hyperlinks.Item(2).TextToDisplay = "f84312db-1505-4f53-a880-5d7ff6201955.3c3df32c68054f28c47665f5d132b728da9dbe5c.zh-cn.xlf"
$ws2.Hyperlinks.Item(3).TextToDisplay = "f84312db-1505-4f53-a880-5d7ff6201955.md"
$ws2.Hyperlinks.Item(4).TextToDisplay = "f84312db-1505-4f53-a880-5d7ff6201955.3c3df32c68054f28c47665f5d132b728da9dbe5c.zh-cn.xlf"

$ws2.Hyperlinks.Item(5).TextToDisplay = "b5a67361-cacf-47e3-aedd-3e6e3b8570f1.md"
$ws2.Hyperlinks.Item(6).TextToDisplay = "b5a67361-cacf-47e3-aedd-3e6e3b8570f1.a76c1a906125cd2ada864bd513b0efbeccfc2b2a.zh-cn.xlf"
$ws2.Hyperlinks.Item(7).TextToDisplay = "b5a67361-cacf-47e3-aedd-3e6e3b8570f1.md"
$ws2.Hyperlinks.Item(8).TextToDisplay = "b5a67361-cacf-47e3-aedd-3e6e3b8570f1.a76c1a906125cd2ada864bd513b0efbeccfc2b2a.zh-cn.xlf"

# Non-hyperlinked cells for row 2 (now f84312db) / row 3 (now b5a67361).
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("E2").Value = "2016-03-23 12:51:09"
$ws2.Range("H2").Value = "2016-03-23 12:51:47"

$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("E3").Value = "2016-03-23 12:52:45"
$ws2.Range("H3").Value = "2016-03-23 12:51:47"

# ---------------- de-de sheet ----------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Item(1).TextToDisplay = "f84312db-1505-4f53-a880-5d7ff6201955.md"
$ws3.Hyperlinks.Item(2).TextToDisplay = "f84312db-1505-4f53-a880-5d7ff6201955.3c3df32c68054f28c47665f5d132b728da9dbe5c.de-de.xlf"
$ws3.Hyperlinks.Item(3).TextToDisplay = "f84312db-1505-4f53-a880-5d7ff6201955.md"
$ws3.Hyperlinks.Item(4).TextToDisplay = "f84312db-1505-4f53-a880-5d7ff6201955.3c3df32c68054f28c47665f5d132b728da9dbe5c.de-de.xlf"

$ws3.Hyperlinks.Item(5).TextToDisplay = "b5a67361-cacf-47e3-aedd-3e6e3b8570f1.md"
$ws3.Hyperlinks.Item(6).TextToDisplay = "b5a67361-cacf-47e3-aedd-3e6e3b8570f1.a76c1a906125cd2ada864bd513b0efbeccfc2b2a.de-de.xlf"
$ws3.Hyperlinks.Item(7).TextToDisplay = "b5a67361-cacf-47e3-aedd-3e6e3b8570f1.md"
$ws3.Hyperlinks.Item(8).TextToDisplay = "b5a67361-cacf-47e3-aedd-3e6e3b8570f1.a76c1a906125cd2ada864bd513b0efbeccfc2b2a.de-de.xlf"

$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("E2").Value = "2016-03-23 12:51:13"
$ws3.Range("H2").Value = "2016-03-23 12:51:56"

$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("E3").Value = "2016-03-23 12:52:49"
$ws3.Range("H3").Value = "2016-03-23 12:51:56"
